# Rename the mis-named "Tools" sheet to "Tool"
$wb = $excel.ActiveWorkbook

$toolSheet = $wb.Worksheets.Item("Tools")
$toolSheet.Name = "Tool"

# Select D14 on the Tool sheet, removing the stale top-left freeze/scroll
$toolSheet.Range("D14").Select()

# Make Tool the active sheet/tab (was previously Purchase)
$toolSheet.Activate()
